$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) The previous update batch (rows 576:589) gets its "Fecha" timestamp
#    value refreshed to a (almost identical) recomputed value.
# ---------------------------------------------------------------------------
$oldBatchValue = 44232.53506424768
for ($r = 576; $r -le 589; $r++) {
    $ws.Cells.Item($r, 4).Value2 = $oldBatchValue
}

# ---------------------------------------------------------------------------
# 2) Append a brand-new update batch: 14 rows (590:603) following the same
#    Nombre / URL / Disponibilidad / Fecha cycle used throughout the sheet.
# ---------------------------------------------------------------------------
$names = @(
    "Odoo",
    "Blackbox",
    "PowerBI",
    "Dropbox",
    "Odoo",
    "GEE",
    "UtilidadesOdoo",
    "Filtros Dashboard",
    "MapStore",
    "GeoServer",
    "Tomcat",
    "Shiny",
    "Github",
    "EZ Exporter"
)

$displayUrls = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

# Hyperlink target address (without any fragment) for every row; the
# MapStore row (index 8) keeps its target without the trailing "#/" and
# uses a separate SubAddress ("/") instead, same as the rest of the sheet.
$linkAddresses = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

$newBatchValue = 44232.55611753627
$firstNewRow = 590

for ($i = 0; $i -lt 14; $i++) {
    $row = $firstNewRow + $i

    $ws.Cells.Item($row, 1).Value2 = $names[$i]
    $ws.Cells.Item($row, 3).Value2 = "Disponible"

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value2 = $newBatchValue
    $dCell.NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value2 = $displayUrls[$i]
    if ($i -eq 8) {
        $ws.Hyperlinks.Add($bCell, $linkAddresses[$i], "/")
    } else {
        $ws.Hyperlinks.Add($bCell, $linkAddresses[$i])
    }
    $bCell.Style = "Hyperlink"
}
